$d = $word.ActiveDocument

# --- Change 1: wrap "spaghetti.f" in curly single-quotes -------------------
$old1 = " that was implemented in spaghetti.f "
$new1 = " that was implemented in ‘spaghetti.f’ "
$r1 = $d.Content
$found1 = $r1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Output "Change1 found: $found1"

# --- Change 2: merge trailing " See below:" into previous sentence ---------
# (no visible text change -- purely a run merge, so nothing required here,
#  but run an execute to normalise formatting marks anyway.)
$old2 = "correctly sorted data. See below:"
$new2 = "correctly sorted data. See below:"
$r2 = $d.Content
$found2 = $r2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Output "Change2 found: $found2"

# --- Change 3: expand the "limitations" paragraph ---------------------------
$old3 = "The limitations of this code is obviously the readability. It took several hours of parsing through the code in order to discover what the algorithm was doing and it was an algorithm that I have been familiar with now for several years."
$new3 = " `tOne of the limitations of this code is obviously the readability. It took several hours of parsing through the code in order to discover what the algorithm was doing and it was an algorithm that I have been familiar with now for several years. Additionally, there was a limitation on the domain of input; only 3 character values can be input into the program or else they are cut off before they were interpreted. This means that domain of values that were accepted ranged from negative 99 through positive 999. Values that included any characters or decimals caused an error in the program and it’s immediate termination. "
$r3 = $d.Content
$found3 = $r3.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
Write-Output "Change3 found: $found3"
